$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update BOM table cell values to reflect reordered / new components
$ws.Range("B1").Value2 = "Bauteil"
$ws.Range("C1").Value2 = "Wert"
$ws.Range("D1").Value2 = "Positionen"
$ws.Range("E1").Value2 = "Reichelt Best Nr."
$ws.Range("F1").Value2 = "Conrad Best Nr."
$ws.Range("G1").Value2 = "Geprüft"
$ws.Range("A2").Value2 = 1
$ws.Range("B2").Value2 = "LDO"
$ws.Range("C2").Value2 = "3.3V 1A"
$ws.Range("D2").Value2 = "U2"
$ws.Range("E2").Value2 = "ZLDO1117G33TA"
$ws.Range("G2").Value2 = "ok"
$ws.Range("A3").Value2 = 3
$ws.Range("B3").Value2 = "C"
$ws.Range("C3").Value2 = "4,7u"
$ws.Range("D3").Value2 = "C1, C2, C51"
$ws.Range("E3").Value2 = "KEM X5R0805 4,7U"
$ws.Range("G3").Value2 = "ok"
$ws.Range("A4").Value2 = 1
$ws.Range("B4").Value2 = "StepDown"
$ws.Range("C4").Value2 = "5V 2A"
$ws.Range("D4").Value2 = "U3"
$ws.Range("E4").Value2 = "TS2596SCS50"
$ws.Range("G4").Value2 = "ok"
$ws.Range("A5").Value2 = 2
$ws.Range("B5").Value2 = "D"
$ws.Range("C5").Value2 = "40V 3A Schottky"
$ws.Range("D5").Value2 = "D1, D44"
$ws.Range("E5").Value2 = "B340A-13-F DII"
$ws.Range("G5").Value2 = "ok"
$ws.Range("A6").Value2 = 2
$ws.Range("B6").Value2 = "C"
$ws.Range("C6").Value2 = "470u"
$ws.Range("D6").Value2 = "C3, C5"
$ws.Range("E6").Value2 = "FK 470/6,3 SP"
$ws.Range("G6").Value2 = "ok"
$ws.Range("A7").Value2 = 5
$ws.Range("C7").Value2 = "100p"
$ws.Range("D7").Value2 = "C4, C6"
$ws.Range("E7").Value2 = "FK 470/6,3 SP"
$ws.Range("G7").Value2 = "ok"
$ws.Range("A8").Value2 = 42
$ws.Range("B8").Value2 = "R"
$ws.Range("C8").Value2 = "91R"
$ws.Range("D8").Value2 = "R1, …, R42"
$ws.Range("E8").Value2 = "RND 155HP05 EQ"
$ws.Range("G8").Value2 = "ok"
$ws.Range("A9").Value2 = 44
$ws.Range("B9").Value2 = "C"
$ws.Range("C9").Value2 = "100n"
$ws.Range("D9").Value2 = "C7, …, C50"
$ws.Range("E9").Value2 = "KEM X7R0805 100N"
$ws.Range("G9").Value2 = "ok"
$ws.Range("A10").Value2 = 1
$ws.Range("B10").Value2 = "U"
$ws.Range("C10").Value2 = "FT232RL"
$ws.Range("D10").Value2 = "U4"
$ws.Range("E10").Value2 = "FT 232 RL"
$ws.Range("G10").Value2 = "ok"
$ws.Range("A11").Value2 = 1
$ws.Range("B11").Value2 = "J"
$ws.Range("C11").Value2 = "USB Type B"
$ws.Range("D11").Value2 = "J3"
$ws.Range("E11").Value2 = "RND 205-00858"
$ws.Range("G11").Value2 = "ok"
$ws.Range("A12").Value2 = 2
$ws.Range("B12").Value2 = "D"
$ws.Range("C12").Value2 = "LED Grün"
$ws.Range("D12").Value2 = "D45, D47"
$ws.Range("E12").Value2 = "EVL 17-21SYGC/S2"
$ws.Range("G12").Value2 = "ok"
$ws.Range("A13").Value2 = 1
$ws.Range("B13").Value2 = "D"
$ws.Range("C13").Value2 = "LED Rot"
$ws.Range("D13").Value2 = "D46"
$ws.Range("E13").Value2 = "KBT KP-2012EC"
$ws.Range("G13").Value2 = "ok"
$ws.Range("A14").Value2 = 3
$ws.Range("B14").Value2 = "R"
$ws.Range("C14").Value2 = "300R"
$ws.Range("D14").Value2 = "R43, R44, R49"
$ws.Range("E14").Value2 = "RND 0805 5 300"
$ws.Range("G14").Value2 = "ok"
$ws.Range("A15").Value2 = 3
$ws.Range("B15").Value2 = "SW"
$ws.Range("C15").Value2 = "SMD Taster"
$ws.Range("D15").Value2 = "SW1 - SW3"
$ws.Range("E15").Value2 = "TASTER 1612.11"
$ws.Range("G15").Value2 = "ok"
$ws.Range("A16").Value2 = 2
$ws.Range("B16").Value2 = "Q"
$ws.Range("C16").Value2 = "NPN Transistor"
$ws.Range("D16").Value2 = "Q1, Q2"
$ws.Range("E16").Value2 = "BCX 19 SMD"
$ws.Range("G16").Value2 = "ok"
$ws.Range("A17").Value2 = 6
$ws.Range("B17").Value2 = "R"
$ws.Range("C17").Value2 = "10k"
$ws.Range("D17").Value2 = "R45 - R48, R50, R51"
$ws.Range("E17").Value2 = "RND 1550805 DN"
$ws.Range("G17").Value2 = "ok"
$ws.Range("A18").Value2 = 1
$ws.Range("B18").Value2 = "L"
$ws.Range("C18").Value2 = "33uH"
$ws.Range("D18").Value2 = "L1"
$ws.Range("E18").Value2 = "L-PIHV4119 33µ"
$ws.Range("G18").Value2 = "ok"
$ws.Range("A19").Value2 = 1
$ws.Range("B19").Value2 = "J"
$ws.Range("C19").Value2 = "Hohlbuchse"
$ws.Range("D19").Value2 = "J1"
$ws.Range("E19").Value2 = "DC-BU 072759"
$ws.Range("G19").Value2 = "ok"
$ws.Range("A20").Value2 = 1
$ws.Range("B20").Value2 = "J"
$ws.Range("C20").Value2 = "Anschlussklemme"
$ws.Range("D20").Value2 = "J2"
$ws.Range("E20").Value2 = "AKL 057-02"
$ws.Range("G20").Value2 = "ok"
$ws.Range("A21").Value2 = 1
$ws.Range("B21").Value2 = "SW"
$ws.Range("C21").Value2 = "THT Schalter 2 Pol"
$ws.Range("D21").Value2 = "SW4"
$ws.Range("E21").Value2 = "SS ESP20"
$ws.Range("G21").Value2 = "ok"
$ws.Range("A22").Value2 = 1
$ws.Range("B22").Value2 = "U"
$ws.Range("C22").Value2 = "OP07"
$ws.Range("D22").Value2 = "U5"
$ws.Range("E22").Value2 = "OP 07 CD"
$ws.Range("G22").Value2 = "ok"
$ws.Range("A23").Value2 = 1
$ws.Range("B23").Value2 = "Q"
$ws.Range("C23").Value2 = "P-FET"
$ws.Range("D23").Value2 = "Q3"
$ws.Range("E23").Value2 = "NDS 0610 SMD"
$ws.Range("G23").Value2 = "ok"
$ws.Range("A24").Value2 = 42
$ws.Range("B24").Value2 = "D"
$ws.Range("C24").Value2 = "RGB"
$ws.Range("D24").Value2 = "D2, …, D43"
$ws.Range("F24").Value2 = "800167562 - VQ"
$ws.Range("G24").Value2 = "ok"
$ws.Range("A25").Value2 = 1
$ws.Range("B25").Value2 = " µC"
$ws.Range("C25").Value2 = "ESP32 WROOM 32D"
$ws.Range("D25").Value2 = "U1"
$ws.Range("F25").Value2 = "1925467 - VQ"
$ws.Range("G25").Value2 = "ok"

# Clear cells that no longer hold data after the reshuffle
$ws.Range("B7").Value2 = $null
$ws.Range("F4").Value2 = $null
$ws.Range("F9").Value2 = $null

# Resize the BOM table to include the two newly added rows (OP07, P-FET)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G25"))

# Adjust column D width slightly (cosmetic resize observed after edit)
$ws.Columns.Item(4).ColumnWidth = 17.9

# Restore the active cell selection left by the editor
$ws.Range("E7").Select()
